$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet references
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)   # 01_To_10_Dec-2020
$ws2 = $wb.Worksheets.Item(2)   # 11_To_20_Dec-2020

# ---------------------------------------------------------------------------
# New order rows for 12 Dec 2020, appended after the existing 16 rows
# (row 17 is the last existing data row) on the "11_To_20_Dec-2020" sheet.
# Columns: A S.No. | B Order ID | C Retailer Name | D Order Amount |
#          E Date/Time | F FOS Name | G Type
# ---------------------------------------------------------------------------
$hyperlinkTarget = "https://fiori.jioconnect.com/sap/bc/ui5_ui5/sap/zehys_dashboard/javascript:void(0);"

$newRows = @(
    @{ Row=18; No=17; OrderId=2004246801; Retailer="CHANDAN TELECOM-(661644693)";               Amount=1040; DateTime="12 Dec 2020 / 12:21"; Fos="";               Type="Manual" },
    @{ Row=19; No=18; OrderId=2004307224; Retailer="chintu electronics-(662114680)";             Amount=1040; DateTime="12 Dec 2020 / 15:17"; Fos="";               Type="Auto" },
    @{ Row=20; No=19; OrderId=2004377153; Retailer="ANISH MOBILE CENTER-(661474493)";            Amount=2080; DateTime="12 Dec 2020 / 15:28"; Fos="KUMAR RAJESH";     Type="Auto" },
    @{ Row=21; No=20; OrderId=2004380406; Retailer="SINGH COMMUNICATION-(660315658)";            Amount=5200; DateTime="12 Dec 2020 / 15:29"; Fos="KUMAR RAJESH";     Type="Auto" },
    @{ Row=22; No=21; OrderId=2004381282; Retailer="ROHIT SHOE STORE-(661670729)";               Amount=5200; DateTime="12 Dec 2020 / 15:30"; Fos="Kumar Jitendra";   Type="Auto" },
    @{ Row=23; No=22; OrderId=2004383653; Retailer="MAA TELECOM-(661361714)";                    Amount=2080; DateTime="12 Dec 2020 / 15:31"; Fos="Kumar Jitendra";   Type="Auto" },
    @{ Row=24; No=23; OrderId=2004383737; Retailer="UTKARSH ENTERPRISES-(661673664)";            Amount=2080; DateTime="12 Dec 2020 / 15:31"; Fos="KUMAR RAJESH";     Type="Auto" },
    @{ Row=25; No=24; OrderId=2004385244; Retailer="KALLU JI DOWN LOADING CENTRE-(661066808)";   Amount=5200; DateTime="12 Dec 2020 / 15:32"; Fos="Kumar Jitendra";   Type="Auto" },
    @{ Row=26; No=25; OrderId=2004393199; Retailer="SUMAN PHONE GHAR-(660315655)";               Amount=5200; DateTime="12 Dec 2020 / 15:36"; Fos="KUMAR RAJESH";     Type="Auto" },
    @{ Row=27; No=26; OrderId=2004449754; Retailer="MEDIA MOBILE GALLERY-(660315675)";           Amount=3120; DateTime="12 Dec 2020 / 17:17"; Fos="KUMAR RAJESH";     Type="Auto" },
    @{ Row=28; No=27; OrderId=2004454152; Retailer="ASHU COMMUNICATION-(661025161)";             Amount=1040; DateTime="12 Dec 2020 / 17:19"; Fos="KUMAR RAJESH";     Type="Auto" },
    @{ Row=29; No=28; OrderId=2004458040; Retailer="SHASHI PAY PHONE-(661176231)";               Amount=5200; DateTime="12 Dec 2020 / 17:21"; Fos="Kumar Jitendra";   Type="Auto" },
    @{ Row=30; No=29; OrderId=2004458444; Retailer="SUBODH PAN-(661562286)";                     Amount=2080; DateTime="12 Dec 2020 / 17:21"; Fos="Kumar Jitendra";   Type="Auto" },
    @{ Row=31; No=30; OrderId=2004461115; Retailer="PRAKASH GENERAL SRINGAR STORE-(661303685)";  Amount=3120; DateTime="12 Dec 2020 / 17:23"; Fos="KUMAR RAJESH";     Type="Auto" },
    @{ Row=32; No=31; OrderId=2004461755; Retailer="chintu electronics-(662114680)";             Amount=1040; DateTime="12 Dec 2020 / 17:23"; Fos="";               Type="Auto" }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $cellA = $ws2.Cells.Item($row, 1)
    $cellA.Value = $r.No
    $cellA.HorizontalAlignment = -4108
    $cellA.VerticalAlignment = -4108

    $cellB = $ws2.Cells.Item($row, 2)
    $cellB.Value = $r.OrderId
    $cellB.HorizontalAlignment = -4108
    $cellB.VerticalAlignment = -4108

    $cellC = $ws2.Cells.Item($row, 3)
    $cellC.Value = $r.Retailer
    $cellC.HorizontalAlignment = -4108
    $cellC.VerticalAlignment = -4108

    $cellD = $ws2.Cells.Item($row, 4)
    $cellD.Value = $r.Amount
    $cellD.HorizontalAlignment = -4108
    $cellD.VerticalAlignment = -4108

    $cellE = $ws2.Cells.Item($row, 5)
    $cellE.Value = $r.DateTime
    $cellE.HorizontalAlignment = -4108
    $cellE.VerticalAlignment = -4108

    $cellF = $ws2.Cells.Item($row, 6)
    if ($r.Fos -ne "") {
        $cellF.Value = $r.Fos
    }
    $cellF.HorizontalAlignment = -4108
    $cellF.VerticalAlignment = -4108

    $cellG = $ws2.Cells.Item($row, 7)
    $cellG.Value = $r.Type
    $cellG.HorizontalAlignment = -4108
    $cellG.VerticalAlignment = -4108

    # Matching hyperlink on the Order ID cell, same as every prior row.
    $ws2.Hyperlinks.Add($cellB, $hyperlinkTarget, [Type]::Missing, $hyperlinkTarget, [string]$r.OrderId)
}

# Trailing row with only the running S.No. counter (no order that day).
$cellA33 = $ws2.Cells.Item(33, 1)
$cellA33.Value = 32
$cellA33.HorizontalAlignment = -4108
$cellA33.VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# View state: "01_To_10_Dec-2020" becomes the active sheet with its header
# row frozen and E1 selected; "11_To_20_Dec-2020" keeps I12 selected.
# ---------------------------------------------------------------------------
$ws2.Activate()
$ws2.Range("I12").Select()

$ws1.Activate()
$ws1.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws1.Range("E1").Select()
